# "Add files via upload" - adds new machine-learning video-tutorial rows
# (topic + Google-Drive hyperlink) to Sheet5, plus an author credit and
# some column sizing / selection bookkeeping that goes with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet5")

# New rows under the existing header (row 4: Topic | GD link | Youtube link)
$ws.Range("A5").Value = "INTRODUCTION TO MACHINE LEARNING"
$ws.Range("C5").Value = "https://drive.google.com/drive/folders/1cd-J8_d6fwlmoQyYsPrYruewp4zYegj0"

$ws.Range("A6").Value = "WHAT IS DATA"
$ws.Range("C6").Value = "https://drive.google.com/drive/folders/1cd-J8_d6fwlmoQyYsPrYruewp4zYegj0"

$ws.Range("A7").Value = "WHAT IS MEANING OF MACHINE LEARNING"
$ws.Range("C7").Value = "https://drive.google.com/drive/u/0/folders/1cd-J8_d6fwlmoQyYsPrYruewp4zYegj0"

$ws.Range("A8").Value = "WHAT IS SUPERVISED LEARNING"
$ws.Range("C8").Value = "https://drive.google.com/drive/u/0/folders/1cd-J8_d6fwlmoQyYsPrYruewp4zYegj0"

# Author credit just under the sheet's CSS title
$ws.Range("A2").Value = "SHREYA"

# Turn the GD-link cells into real hyperlinks (applies the built-in
# "Hyperlink" style too)
$ws.Hyperlinks.Add($ws.Range("C6"), "https://drive.google.com/drive/folders/1cd-J8_d6fwlmoQyYsPrYruewp4zYegj0") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "https://drive.google.com/drive/u/0/folders/1cd-J8_d6fwlmoQyYsPrYruewp4zYegj0") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://drive.google.com/drive/u/0/folders/1cd-J8_d6fwlmoQyYsPrYruewp4zYegj0") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://drive.google.com/drive/folders/1cd-J8_d6fwlmoQyYsPrYruewp4zYegj0") | Out-Null

# Widen the Topic / link columns so the new long text/links are readable
$ws.Columns.Item(2).ColumnWidth = 28.109375
$ws.Columns.Item(3).ColumnWidth = 70.21875

# Leave the selection where the author ended up after entering the data
$ws.Range("E10").Select() | Out-Null
